# Update cryptocurrency price/volume data to reflect the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.807.70"
$ws.Range("E2").Value = "  -0.33%  "
$ws.Range("D3").Value = "1.635.52"
$ws.Range("E3").Value = "  +0.11%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'215.38"
$ws.Range("E5").Value = "  -0.31%  "
$ws.Range("D6").Value = "'0.5049"
$ws.Range("E6").Value = "  -0.39%  "
$ws.Range("D7").Value = "'1.004"
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("D8").Value = "'0.2578"
$ws.Range("E8").Value = "  +0.19%  "
$ws.Range("D9").Value = "'0.06414"
$ws.Range("E9").Value = "  +0.92%  "
$ws.Range("D10").Value = "'20.34"
$ws.Range("E10").Value = "  +3.90%  "
$ws.Range("D11").Value = "'0.07803"
$ws.Range("E11").Value = "  +0.65%  "
$ws.Range("D12").Value = "'4.293"
$ws.Range("E12").Value = "  +0.98%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "1.860.64"
$ws.Range("E13").Value = "  +0.00%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.634.30"
$ws.Range("E14").Value = "  -0.31%  "
$ws.Range("D15").Value = "'0.5633"
$ws.Range("E15").Value = "  +2.41%  "
$ws.Range("D16").Value = "0.0₅7630"
$ws.Range("E16").Value = "  -0.87%  "
$ws.Range("D17").Value = "'63.04"
$ws.Range("E17").Value = "  -1.21%  "
$ws.Range("D18").Value = "25.816.17"
$ws.Range("E18").Value = "  -0.39%  "
$ws.Range("E19").Value = "  +0.13%  "
$ws.Range("D20").Value = "'194.46"
$ws.Range("E20").Value = "  +0.18%  "
$ws.Range("D21").Value = "'4.346"
$ws.Range("E21").Value = "  -1.93%  "
$ws.Range("D22").Value = "'9.919"
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("D23").Value = "'6.103"
$ws.Range("E23").Value = "  +0.92%  "
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("D25").Value = "'1.789"
$ws.Range("E25").Value = "  -6.08%  "
$ws.Range("D26").Value = "'140.39"
$ws.Range("E26").Value = "  -1.52%  "
$ws.Range("D27").Value = "'0.1254"
$ws.Range("E27").Value = "  +1.29%  "
$ws.Range("D28").Value = "'6.824"
$ws.Range("E28").Value = "  +0.29%  "
$ws.Range("D29").Value = "'15.43"
$ws.Range("E29").Value = "  -0.82%  "
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("D31").Value = "'0.04912"
$ws.Range("E31").Value = "  +0.48%  "
$ws.Range("D32").Value = "'3.313"
$ws.Range("E32").Value = "  +1.84%  "
$ws.Range("D33").Value = "'3.248"
$ws.Range("E33").Value = "  +1.87%  "
$ws.Range("D34").Value = "'1.574"
$ws.Range("E34").Value = "  +2.11%  "
$ws.Range("D35").Value = "'2.380"
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("D36").Value = "'0.9049"
$ws.Range("E36").Value = "  +0.29%  "
$ws.Range("D37").Value = "'2.568"
$ws.Range("E37").Value = "  +0.42%  "
$ws.Range("D38").Value = "'0.5536"
$ws.Range("E38").Value = "  +0.73%  "
$ws.Range("D39").Value = "1.125.51"
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("D40").Value = "'0.01559"
$ws.Range("E40").Value = "  +0.33%  "
$ws.Range("D41").Value = "'1.002"
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("E42").Value = "  -1.52%  "
$ws.Range("D43").Value = "'0.8013"
$ws.Range("E43").Value = "  -0.24%  "
$ws.Range("E44").Value = "  +0.66%  "
$ws.Range("E45").Value = "  -0.07%  "
$ws.Range("D46").Value = "0.0₈112"
$ws.Range("E46").Value = "  -6.96%  "
$ws.Range("D47").Value = "'55.52"
$ws.Range("D48").Value = "'0.4262"
$ws.Range("E48").Value = "  -4.29%  "
$ws.Range("D49").Value = "'7.731"
$ws.Range("E49").Value = "  +2.88%  "
$ws.Range("D50").Value = "'0.05021"
$ws.Range("E50").Value = "  -2.43%  "
$ws.Range("E51").Value = "  +0.59%  "
